$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Formula = "=""0.06"""
$ws.Range("C2").Formula = "=""-0.0"""
$ws.Range("D2").Formula = "=""0.04"""
$ws.Range("E2").Formula = "=""0.14"""
$ws.Range("F2").Formula = "=""0.15"""
$ws.Range("G2").Formula = "=""-0.03"""
$ws.Range("B3").Formula = "=""0.12"""
$ws.Range("C3").Formula = "=""0.03"""
$ws.Range("D3").Formula = "=""0.09"""
$ws.Range("E3").Formula = "=""0.19"""
$ws.Range("F3").Formula = "=""0.11"""
$ws.Range("G3").Formula = "=""-0.15"""
$ws.Range("B4").Formula = "=""0.13"""
$ws.Range("C4").Formula = "=""-0.02"""
$ws.Range("D4").Formula = "=""0.13"""
$ws.Range("E4").Formula = "=""0.25**"""
$ws.Range("F4").Formula = "=""0.02"""
$ws.Range("G4").Formula = "=""-0.23*"""
$ws.Range("B5").Formula = "=""0.11"""
$ws.Range("C5").Formula = "=""-0.12"""
$ws.Range("D5").Formula = "=""0.09"""
$ws.Range("E5").Formula = "=""0.17"""
$ws.Range("F5").Formula = "=""-0.03"""
$ws.Range("G5").Formula = "=""-0.15"""
$ws.Range("B6").Formula = "=""0.01"""
$ws.Range("C6").Formula = "=""-0.22*"""
$ws.Range("D6").Formula = "=""-0.07"""
$ws.Range("E6").Formula = "=""0.08"""
$ws.Range("F6").Formula = "=""0.05"""
$ws.Range("G6").Formula = "=""-0.11"""
$ws.Range("B7").Formula = "=""0.07"""
$ws.Range("C7").Formula = "=""-0.23*"""
$ws.Range("D7").Formula = "=""-0.01"""
$ws.Range("E7").Formula = "=""-0.07"""
$ws.Range("F7").Formula = "=""0.02"""
$ws.Range("G7").Formula = "=""-0.06"""
$ws.Range("B8").Formula = "=""0.01"""
$ws.Range("C8").Formula = "=""-0.28**"""
$ws.Range("D8").Formula = "=""-0.07"""
$ws.Range("E8").Formula = "=""-0.01"""
$ws.Range("F8").Formula = "=""-0.09"""
$ws.Range("G8").Formula = "=""-0.07"""
$ws.Range("B9").Formula = "=""-0.07"""
$ws.Range("C9").Formula = "=""-0.32**"""
$ws.Range("D9").Formula = "=""-0.16"""
$ws.Range("E9").Formula = "=""-0.08"""
$ws.Range("F9").Formula = "=""-0.17"""
$ws.Range("G9").Formula = "=""-0.12"""
$ws.Range("B10").Formula = "=""-0.17"""
$ws.Range("C10").Formula = "=""-0.22*"""
$ws.Range("D10").Formula = "=""-0.18"""
$ws.Range("E10").Formula = "=""-0.03"""
$ws.Range("F10").Formula = "=""-0.11"""
$ws.Range("G10").Formula = "=""-0.13"""
$ws.Range("B11").Formula = "=""-0.25**"""
$ws.Range("C11").Formula = "=""-0.28**"""
$ws.Range("D11").Formula = "=""-0.21*"""
$ws.Range("E11").Formula = "=""-0.07"""
$ws.Range("F11").Formula = "=""-0.17"""
$ws.Range("G11").Formula = "=""-0.05"""
$ws.Range("B12").Formula = "=""-0.3**"""
$ws.Range("C12").Formula = "=""-0.32**"""
$ws.Range("D12").Formula = "=""-0.23*"""
$ws.Range("E12").Formula = "=""-0.04"""
$ws.Range("F12").Formula = "=""-0.18"""
$ws.Range("G12").Formula = "=""-0.14"""
$ws.Range("B13").Formula = "=""-0.21"""
$ws.Range("C13").Formula = "=""-0.24*"""
$ws.Range("D13").Formula = "=""-0.12"""
$ws.Range("E13").Formula = "=""-0.06"""
$ws.Range("F13").Formula = "=""-0.28**"""
$ws.Range("G13").Formula = "=""-0.05"""
$srcRange = $ws.Range("B2:G13")
$srcRange.Copy() | Out-Null
$srcRange.PasteSpecial(-4163) | Out-Null
